$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": insert new client row before the summary row ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(11).Insert()

$ws1.Cells.Item(11,1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(11,2).Value = "VEHINVER SA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(11,$c).Value = 0
}

# The old "X de 9" summary row got pushed down to row 12 - bump the count text to "de 10"
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(12,$c)
    $cell.Value = ($cell.Value2 -replace "de 9", "de 10")
}

# --- Sheet "VENTA MENSUAL": insert new client row before the totals row ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(11).Insert()

$ws2.Cells.Item(11,1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(11,2).Value = "VEHINVER SA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(11,$c).Value = 0
}
